$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 592
$ws.Range("F3").Value = 5651
$ws.Range("F5").Value = 459
$ws.Range("F12").Value = 3094
$ws.Range("F13").Value = 1926
$ws.Range("F17").Value = 53
$ws.Range("F18").Value = 155
$ws.Range("F20").Value = 979
$ws.Range("F21").Value = 355
$ws.Range("F23").Value = 21
$ws.Range("F24").Value = 3582
$ws.Range("F25").Value = 1138
$ws.Range("F26").Value = 2839
$ws.Range("F27").Value = 284
$ws.Range("F28").Value = 2207
$ws.Range("F29").Value = 4128
$ws.Range("F31").Value = 922
$ws.Range("F32").Value = 472
$ws.Range("F33").Value = 1314
$ws.Range("F34").Value = 73
$ws.Range("F36").Value = 1012
$ws.Range("F37").Value = 1280
$ws.Range("F39").Value = 1070
$ws.Range("F40").Value = 692
$ws.Range("F41").Value = 570
$ws.Range("F42").Value = 419
$ws.Range("F43").Value = 6
$ws.Range("F44").Value = 63
$ws.Range("F45").Value = 319
$ws.Range("F46").Value = 3582

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 5
$ws.Range("F10").Value = 908
$ws.Range("F20").Value = 2
$ws.Range("F25").Value = 13

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 592
$ws.Range("F3").Value = 5652
$ws.Range("F10").Value = 3094
$ws.Range("F12").Value = 1926
$ws.Range("F16").Value = 908
$ws.Range("F18").Value = 155
$ws.Range("F19").Value = 979
$ws.Range("F20").Value = 355
$ws.Range("F21").Value = 3582
$ws.Range("F24").Value = 1138
$ws.Range("F26").Value = 2839
$ws.Range("F27").Value = 2207
$ws.Range("F28").Value = 4128
$ws.Range("F31").Value = 922
$ws.Range("F32").Value = 1314
$ws.Range("F34").Value = 1012
$ws.Range("F35").Value = 1280
$ws.Range("F37").Value = 1070
$ws.Range("F39").Value = 692
$ws.Range("F41").Value = 419
$ws.Range("F44").Value = 13
$ws.Range("F45").Value = 63
$ws.Range("F47").Value = 319
$ws.Range("F48").Value = 3582
